$d = $word.ActiveDocument

$replacements = @(
    @{old = "881÷5=176, 1"; new = "427÷3=142, 1"},
    @{old = "559÷9=62, 1"; new = "812÷3=270, 2"},
    @{old = "592÷9=65, 7"; new = "733÷4=183, 1"},
    @{old = "822÷6=137, 0"; new = "359÷2=179, 1"},
    @{old = "633÷2=316, 1"; new = "540÷8=67, 4"},
    @{old = "993÷9=110, 3"; new = "820÷4=205, 0"},
    @{old = "961÷8=120, 1"; new = "541÷3=180, 1"},
    @{old = "908÷4=227, 0"; new = "846÷7=120, 6"},
    @{old = "249÷6=41, 3"; new = "742÷6=123, 4"},
    @{old = "811÷9=90, 1"; new = "835÷8=104, 3"},
    @{old = "873÷4=218, 1"; new = "416÷8=52, 0"},
    @{old = "686÷4=171, 2"; new = "416÷3=138, 2"},
    @{old = "560÷6=93, 2"; new = "634÷9=70, 4"},
    @{old = "996÷7=142, 2"; new = "971÷4=242, 3"},
    @{old = "751÷7=107, 2"; new = "418÷5=83, 3"},
    @{old = "622÷4=155, 2"; new = "938÷5=187, 3"},
    @{old = "316÷3=105, 1"; new = "897÷6=149, 3"},
    @{old = "640÷8=80, 0"; new = "313÷7=44, 5"},
    @{old = "575÷2=287, 1"; new = "517÷9=57, 4"},
    @{old = "445÷9=49, 4"; new = "418÷3=139, 1"},
    @{old = "898÷3=299, 1"; new = "587÷7=83, 6"},
    @{old = "888÷9=98, 6"; new = "514÷3=171, 1"},
    @{old = "698÷7=99, 5"; new = "842÷2=421, 0"},
    @{old = "623÷6=103, 5"; new = "186÷8=23, 2"},
    @{old = "486÷5=97, 1"; new = "606÷3=202, 0"}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
